$wb = $excel.ActiveWorkbook

$wsFuzzy = $wb.Worksheets.Item("Fuzzificación")
$wsSistema = $wb.Worksheets.Item("Sistema Difuso")

# --- Fuzzificación: add the "Delicioso" membership formula in column H ---

# H2 is filled on its own (mirrors the existing single-cell formulas already in B2/C2/D2/G2)
$wsFuzzy.Range("H2").Formula = "=IF(F2<7,0,IF(F2<9,(F2-7)/2,IF(F2<=10,1,0)))"

# Give H3:H102 the same visual style already used by H2 (blue-filled cell, style index 3)
# before filling them with formulas, so they pick up that format.
$wsFuzzy.Range("H2").Copy()
$wsFuzzy.Range("H3:H102").PasteSpecial(-4122)  # xlPasteFormats

# H3:H66 filled as one block (matches the first fill-down boundary already used by columns B/C/D/G)
$wsFuzzy.Range("H3:H66").Formula = "=IF(F3<7,0,IF(F3<9,(F3-7)/2,IF(F3<=10,1,0)))"

# H67:H102 filled as the second block (matches the second fill-down boundary already used by columns B/C/D/G)
$wsFuzzy.Range("H67:H102").Formula = "=IF(F67<7,0,IF(F67<9,(F67-7)/2,IF(F67<=10,1,0)))"

# Move the selection on "Fuzzificación" from G2 to B2
$wsFuzzy.Range("B2").Select()

# Make "Sistema Difuso" the active sheet/tab, with E3 selected
$wsSistema.Activate()
$wsSistema.Range("E3").Select()
